$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 135 (shifts existing rows 135-202 down to 136-203)
$ws.Rows.Item(135).Insert()

# Populate the new row 135 with the new weekly record
$ws.Range("A135").Value = 8
$ws.Range("B135").Value = "Terminal La Palmera de La Serena"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 44572
$ws.Range("E135").Value = 4
$ws.Range("F135").Value = 100112003
$ws.Range("G135").Value = "Ajo"
$ws.Range("H135").Value = "Chino"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 540
$ws.Range("K135").Value = 18000
$ws.Range("L135").Value = 19000
$ws.Range("M135").Value = 18500
$ws.Range("N135").Value = "$/caja 10 kilos"
$ws.Range("O135").Value = "China"
$ws.Range("P135").Value = 1850
$ws.Range("Q135").Value = 10
$ws.Range("R135").Value = "Hortaliza"
